# Auto-update gym prices
#
# Updates 4 price cells across two sheets:
#   "4x4 Squat Racks"!C2  $2,155.00 -> $2,153.00
#   "4x4 Squat Racks"!C4  "Price not available" -> $1,520.00
#   "Squat Stands"!C2     $1,558.00 -> $1,557.00
#   "Squat Stands"!C6     "Price not available" -> $820.00
#
# The price column stores plain text (e.g. "$2,153.00"), not numeric
# currency values, so each cell's number format is forced to Text ("@")
# before the value is assigned. Without this, Excel's automatic
# number-detection would silently reinterpret a "$"-and-comma string as a
# numeric currency value instead of keeping the literal text.

$wb = $excel.ActiveWorkbook

$wsRacks = $wb.Worksheets.Item("4x4 Squat Racks")
$wsRacks.Range("C2").NumberFormat = "@"
$wsRacks.Range("C2").Value = "$2,153.00"

$wsRacks.Range("C4").NumberFormat = "@"
$wsRacks.Range("C4").Value = "$1,520.00"

$wsStands = $wb.Worksheets.Item("Squat Stands")
$wsStands.Range("C2").NumberFormat = "@"
$wsStands.Range("C2").Value = "$1,557.00"

$wsStands.Range("C6").NumberFormat = "@"
$wsStands.Range("C6").Value = "$820.00"
